# "After writing the second article" -- append the 2024 data block (10 regions x
# 4 fuel types = 40 rows) to Sheet1, continuing directly after the existing
# 241 data rows (header row 1 + rows 2..241).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$regions = @("Baku", "Ganja", "Sumgait", "Mingachevir", "Nakhchivan", "Shaki", "Lankaran", "Yevlakh", "Gabala", "Goychay")
$fuels   = @("Diesel", "Gasoline", "Electric", "Hybrid")

$startRow = 242
$startSeq = 241
$endRow = $startRow + ($regions.Count * $fuels.Count) - 1

# Column A holds the row counter as text (matches the existing rows, which are
# stored as shared strings rather than numbers). Force text storage for the
# whole new block up front, write the values, then restore the default style
# so no cell ends up with a lingering custom format.
$colA = $ws.Range("A$startRow`:A$endRow")
$colA.NumberFormat = "@"

$r = $startRow
$seq = $startSeq
foreach ($region in $regions) {
    foreach ($fuel in $fuels) {
        $ws.Cells.Item($r, 1).Value = "$seq"
        $ws.Cells.Item($r, 2).Value = 2024
        $ws.Cells.Item($r, 3).Value = "`"$region`""
        $ws.Cells.Item($r, 4).Value = "`"$fuel`""
        $r = $r + 1
        $seq = $seq + 1
    }
}

$colA.Style = "Normal"
